$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 1
$ws.Range("C1").Value = 2
$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 4
$ws.Range("F1").Value = 5
$ws.Range("G1").Value = 6
$ws.Range("H1").Value = 7
$ws.Range("I1").Value = 8
$ws.Range("J1").Value = 9
